$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was added to the dataset. It belongs right above the
# current row 155, so insert a fresh row there (this pushes the existing
# rows 155-168 down to 156-169, preserving all their original values).
$ws.Rows(155).Insert()

# Populate the newly inserted row 155 with the new Fruta/Mango record.
$ws.Range("A155").Value = 7
$ws.Range("B155").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C155").Value = "Ñuble"
$ws.Range("D155").Value = 45132
$ws.Range("E155").Value = 16
$ws.Range("F155").Value = "Fruta"
$ws.Range("G155").Value = 100108
$ws.Range("H155").Value = "Tropicales y subtropicales"
$ws.Range("I155").Value = 100108002
$ws.Range("J155").Value = "Mango"
$ws.Range("K155").Value = "Sin especificar"
$ws.Range("L155").Value = "Primera"
$ws.Range("M155").Value = 30
$ws.Range("N155").Value = 8000
$ws.Range("O155").Value = 8000
$ws.Range("P155").Value = 8000
$ws.Range("Q155").Value = "`$/bandeja 4 kilos"
$ws.Range("R155").Value = "Brasil"
$ws.Range("S155").Value = 2000
$ws.Range("T155").Value = 4
